# Updated cryptos list (price + 1h volume columns) per the GitHub Actions refresh.
# For Price cells whose new text looks like a plain decimal number (e.g. "578.43"),
# we briefly force text format so Excel keeps it as a string instead of coercing it
# to a numeric value, then restore the default "Normal" style so no stray
# per-cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.631.62"
$ws.Range("E2").Value = "  +3.98%  "
$ws.Range("D3").Value = "3.257.26"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.79%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.21%  "
$ws.Range("D9").Value = "3.258.29"
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.00%  "
$ws.Range("E11").Value = "  +3.51%  "
$ws.Range("E12").Value = "  +5.00%  "
$ws.Range("D13").Value = "3.833.78"
$ws.Range("E13").Value = "  +3.50%  "
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.56%  "
$ws.Range("D16").Value = "67.618.88"
$ws.Range("E16").Value = "  +4.09%  "
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("D18").Value = "3.262.31"
$ws.Range("E18").Value = "  +3.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.91%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.512"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +6.47%  "
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +5.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.47%  "
$ws.Range("E37").Value = "  +3.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.854"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "365.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.17%  "
$ws.Range("D45").Value = "2.742.69"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0682"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("E50").Value = "  +6.73%  "
$ws.Range("E51").Value = "  +0.31%  "
